$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-13 Tuesday" "2026-01-14 Wednesday"

Replace-Text "22×58=" "34×83="
Replace-Text "64×95=" "86×99="
Replace-Text "56×27=" "68×73="
Replace-Text "57×59=" "78×25="
Replace-Text "84×81=" "87×52="

Replace-Text "30×66=" "78×70="
Replace-Text "41×40=" "37×37="
Replace-Text "80×93=" "35×16="
Replace-Text "71×81=" "59×29="
Replace-Text "40×69=" "49×89="

Replace-Text "99×82=" "47×84="
Replace-Text "31×62=" "48×98="
Replace-Text "29×45=" "14×21="
Replace-Text "58×43=" "22×16="
Replace-Text "43×33=" "30×25="

Replace-Text "48×32=" "38×77="
Replace-Text "36×63=" "47×69="
Replace-Text "78×48=" "75×18="
Replace-Text "30×49=" "43×45="
Replace-Text "71×15=" "54×67="

Replace-Text "68×11=" "97×27="
Replace-Text "15×90=" "78×61="
Replace-Text "13×87=" "12×17="
Replace-Text "67×70=" "69×69="
Replace-Text "35×48=" "90×36="
